# Fill in the Dutch (nl) translations for rows 98-106, column C, which
# previously held the untranslated placeholder ("XXXX" / shared string 0)
# rendered with the red "Bad" style. Also restyle those cells to the
# green "Good"-like look (black text on green fill) now that they are
# translated, and move the sheet's viewport/selection down to that block
# (per the commit: "Dont try to load new content when selecting a new
# location").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$translations = @(
    "Versturen",
    "Verstuurd",
    "Fout opgetreden",
    "Content ophalen...",
    "Controleren op nieuwe content...",
    "Nieuwe content gevonden",
    "Geen nieuwe content gevonden",
    "Optioneel",
    "Voer e-mail in om INSPEC2T updates te ontvangen"
)

$startRow = 98
for ($i = 0; $i -lt $translations.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value = $translations[$i]
}

# Re-style the now-translated cells: black Calibri text on the green fill
# (previously they used the red/black "Bad" highlight reserved for missing
# translations).
$target = $ws.Range("C98:C106")
$target.Font.Name = "Calibri"
$target.Font.Color = 0
$target.Interior.Color = 13561798

# Reflect the editor having scrolled down to / selected this block while
# working on it.
$excel.ActiveWindow.ScrollRow = 88
$excel.ActiveWindow.ScrollColumn = 1
$target.Select()
